$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B13").Value = "21-07-2015"
$ws.Range("C13").Value = "Loc + musique + sfx"
$ws.Range("D13").Value = 4

$ws.Range("D14").Select()
